$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "VIMMP_DEF" in column F, matching the format of E1 (bold/bordered header style)
$ws.Range("F1").Value = "VIMMP_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Fill new column F data rows with "[]" (same plain formatting as column E data cells)
$ws.Range("F2").Value = "[]"
$ws.Range("F3").Value = "[]"
$ws.Range("F4").Value = "[]"
